# adecuar reporte para vacaciones
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asistencia")

$empleado = "Karen Burgos"

# Entrada (col B) / Salida (col C) serial date-times for rows 2..14
# (Tue 2025-07-01 .. Tue 2025-07-15, Sundays skipped, Saturdays half-day)
$filas = 2..14
$entradas = @(
    45839.291666666664,
    45840.291666608799,
    45841.291666608799,
    45842.291666608799,
    45843.291666608799,
    45845.291666666664,
    45846.291666608799,
    45847.291666608799,
    45848.291666608799,
    45849.291666608799,
    45850.291666608799,
    45852.291666666664,
    45853.291666608799
)
$salidas = @(
    45839.666666666664,
    45840.666666608799,
    45841.666666608799,
    45842.666666608799,
    45843.458333333336,
    45845.666666666664,
    45846.666666608799,
    45847.666666608799,
    45848.666666608799,
    45849.666666608799,
    45850.458333333336,
    45852.666666666664,
    45853.666666608799
)

for ($i = 0; $i -lt $filas.Count; $i++) {
    $r = $filas[$i]
    $ws.Cells.Item($r, 1).Value = $empleado
    $ws.Cells.Item($r, 2).Value = $entradas[$i]
    $ws.Cells.Item($r, 3).Value = $salidas[$i]
}

# Bring the Entrada/Salida columns onto the same date/time display format
# (the "yyyy/mm/dd hh:mm:ss" custom format is retired in favour of the
# "d/m/yyyy hh:mm:ss" one already used elsewhere on the sheet).
$ws.Range("B2:C14").NumberFormat = 'd/m/yyyy\ hh:mm:ss'
$wb.DeleteNumberFormat('yyyy/mm/dd\ hh:mm:ss')

# Selection now spans the populated attendance block
$ws.Range("B2:C14").Select()

# The column-E validation (referencing a broken #REF! list) moves out of the
# legacy x14 extension and into the main dataValidations collection.
$dv = $ws.Range("E1:E1048576").Validation
$dv.Delete()
$dv.Add(3, 1, 1, "#REF!")
